$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C (text) - plain value assignment is safe
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# Column D (price) - force text format to avoid numeric auto-coercion
$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '58.753.63'
$cell.Style = "Normal"
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.984.06'
$cell.Style = "Normal"
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '560.89'
$cell.Style = "Normal"
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '137.16'
$cell.Style = "Normal"
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '2.974.25'
$cell.Style = "Normal"
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.131'
$cell.Style = "Normal"
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '4.81'
$cell.Style = "Normal"
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.0000228'
$cell.Style = "Normal"
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '33.60'
$cell.Style = "Normal"
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '3.480.02'
$cell.Style = "Normal"
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '6.95'
$cell.Style = "Normal"
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '2.987.21'
$cell.Style = "Normal"
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '58.831.60'
$cell.Style = "Normal"
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '424.18'
$cell.Style = "Normal"
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '13.52'
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.710'
$cell.Style = "Normal"
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '7.11'
$cell.Style = "Normal"
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '13.39'
$cell.Style = "Normal"
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '80.13'
$cell.Style = "Normal"
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '7.75'
$cell.Style = "Normal"
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '25.65'
$cell.Style = "Normal"
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '6.07'
$cell.Style = "Normal"
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '0.1000'
$cell.Style = "Normal"
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.987'
$cell.Style = "Normal"
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0755'
$cell.Style = "Normal"
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '5.73'
$cell.Style = "Normal"
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '2.07'
$cell.Style = "Normal"
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '48.48'
$cell.Style = "Normal"
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '8.68'
$cell.Style = "Normal"
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '2.75'
$cell.Style = "Normal"
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '400.81'
$cell.Style = "Normal"
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '2.728.11'
$cell.Style = "Normal"
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.0347'
$cell.Style = "Normal"
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '125.57'
$cell.Style = "Normal"
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.242'
$cell.Style = "Normal"
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '2.01'
$cell.Style = "Normal"
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '32.08'
$cell.Style = "Normal"
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '23.23'
$cell.Style = "Normal"

# Column E (volume %) - plain value assignment is safe (leading/trailing spaces + % block numeric parse)
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('E3').Value = '  +3.28%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('E6').Value = '  +11.92%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +4.38%  '
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('E10').Value = '  +5.29%  '
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('E12').Value = '  +3.77%  '
$ws.Range('E13').Value = '  +7.83%  '
$ws.Range('E14').Value = '  +4.29%  '
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('E17').Value = '  +5.90%  '
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('E21').Value = '  +5.00%  '
$ws.Range('E22').Value = '  +6.25%  '
$ws.Range('E23').Value = '  +4.21%  '
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('E25').Value = '  +4.08%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  +8.80%  '
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('E30').Value = '  +7.53%  '
$ws.Range('E31').Value = '  +4.11%  '
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('E34').Value = '  +8.77%  '
$ws.Range('E35').Value = '  +22.03%  '
$ws.Range('E36').Value = '  +6.88%  '
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('E40').Value = '  +15.00%  '
$ws.Range('E41').Value = '  +11.37%  '
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('E43').Value = '  +4.11%  '
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('E45').Value = '  +4.73%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('E47').Value = '  +6.07%  '
$ws.Range('E48').Value = '  +4.44%  '
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('E50').Value = '  +19.94%  '
$ws.Range('E51').Value = '  +2.04%  '

Write-Output "applied 92 cell updates"
